$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (TagId 1): Desc + TagGroupId change, TagName text itself changes ---
# C4 keeps the same shared-string slot but the text of that tag changes to "修炼奇才"
$ws.Cells.Item(4, 3).Value = "修炼奇才"
# D4: Desc "好漂亮的姐姐~" -> "适合修炼"
$ws.Cells.Item(4, 4).Value = "适合修炼"
# E4: TagGroupId 0 -> 1
$ws.Cells.Item(4, 5).Value = 1

# --- New data rows 5-10 (TagId 2-7) ---
# Columns: B=TagId, C=TagName, D=Desc, E=TagGroupId, F=IsGoodTag, G=IsInherit, H=ifNatrual, I=ifPatnerNatrual

# Row 5 - TagId 2
$ws.Cells.Item(5, 2).Value = 2
$ws.Cells.Item(5, 3).Value = "七世好人"
$ws.Cells.Item(5, 4).Value = "灵魂强健"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = $true
$ws.Cells.Item(5, 7).Value = $true
$ws.Cells.Item(5, 8).Value = $true
$ws.Cells.Item(5, 9).Value = $true

# Row 6 - TagId 3
$ws.Cells.Item(6, 2).Value = 3
$ws.Cells.Item(6, 3).Value = "天生灵眼"
$ws.Cells.Item(6, 4).Value = "知识广博"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = $true
$ws.Cells.Item(6, 7).Value = $true
$ws.Cells.Item(6, 8).Value = $true
$ws.Cells.Item(6, 9).Value = $true

# Row 7 - TagId 4
$ws.Cells.Item(7, 2).Value = 4
$ws.Cells.Item(7, 3).Value = "妖灵附体"
$ws.Cells.Item(7, 4).Value = "妖灵亲和"
$ws.Cells.Item(7, 5).Value = 4
$ws.Cells.Item(7, 6).Value = $true
$ws.Cells.Item(7, 7).Value = $true
$ws.Cells.Item(7, 8).Value = $true
$ws.Cells.Item(7, 9).Value = $true

# Row 8 - TagId 5
$ws.Cells.Item(8, 2).Value = 5
$ws.Cells.Item(8, 3).Value = "神仙点化"
$ws.Cells.Item(8, 4).Value = "神力亲和"
$ws.Cells.Item(8, 5).Value = 5
$ws.Cells.Item(8, 6).Value = $true
$ws.Cells.Item(8, 7).Value = $true
$ws.Cells.Item(8, 8).Value = $true
$ws.Cells.Item(8, 9).Value = $true

# Row 9 - TagId 6
$ws.Cells.Item(9, 2).Value = 6
$ws.Cells.Item(9, 3).Value = "九阴之人"
$ws.Cells.Item(9, 4).Value = "鬼魂亲和"
$ws.Cells.Item(9, 5).Value = 6
$ws.Cells.Item(9, 6).Value = $true
$ws.Cells.Item(9, 7).Value = $true
$ws.Cells.Item(9, 8).Value = $true
$ws.Cells.Item(9, 9).Value = $true

# Row 10 - TagId 7
$ws.Cells.Item(10, 2).Value = 7
$ws.Cells.Item(10, 3).Value = "七彩祥云"
$ws.Cells.Item(10, 4).Value = "佛光亲和"
$ws.Cells.Item(10, 5).Value = 7
$ws.Cells.Item(10, 6).Value = $true
$ws.Cells.Item(10, 7).Value = $true
$ws.Cells.Item(10, 8).Value = $true
$ws.Cells.Item(10, 9).Value = $true

# --- Column F width tweak (13.875 -> ~13.887) ---
$ws.Columns.Item(6).ColumnWidth = 13.18

# --- Selection moves to K32 ---
$null = $ws.Range("K32").Select()
